$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.354.12'
$ws.Range('E2').Value = '  -3.06%  '
$ws.Range('D3').Value = '3.164.14'
$ws.Range('E3').Value = '  -4.66%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '571.26'
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('D6').Value = '168.26'
$ws.Range('E6').Value = '  -6.92%  '
$ws.Range('D7').Value = '0.604'
$ws.Range('E7').Value = '  -7.49%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '3.185.13'
$ws.Range('E9').Value = '  -4.02%  '
$ws.Range('E10').Value = '  -5.92%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').Value = '0.389'
$ws.Range('E12').Value = '  -3.23%  '
$ws.Range('D13').Value = '3.712.64'
$ws.Range('E13').Value = '  -4.75%  '
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '64.448.84'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('D16').Value = '25.36'
$ws.Range('E16').Value = '  -4.36%  '
$ws.Range('E17').Value = '  -4.31%  '
$ws.Range('D18').Value = '3.157.59'
$ws.Range('E18').Value = '  -4.84%  '
$ws.Range('D19').Value = '417.36'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = '12.87'
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('E21').Value = '  -3.62%  '
$ws.Range('D22').Value = '7.15'
$ws.Range('E22').Value = '  -3.16%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '69.63'
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('E26').Value = '  -1.02%  '
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('E28').Value = '  -10.39%  '
$ws.Range('D29').Value = '8.81'
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').Value = '1.83'
$ws.Range('E31').Value = '  -4.72%  '
$ws.Range('D32').Value = '21.74'
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('D34').Value = '5.06'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('D35').Value = '6.36'
$ws.Range('E35').Value = '  -4.08%  '
$ws.Range('E36').Value = '  -5.39%  '
$ws.Range('D37').Value = '155.63'
$ws.Range('E37').Value = '  -2.89%  '
$ws.Range('E38').Value = '  -5.23%  '
$ws.Range('D39').Value = '2.702.98'
$ws.Range('E39').Value = '  -5.68%  '
$ws.Range('E40').Value = '  -5.93%  '
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').Value = '23.99'
$ws.Range('E42').Value = '  -9.16%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '39.08'
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.717'
$ws.Range('E44').Value = '  -5.62%  '
$ws.Range('D45').Value = '0.0617'
$ws.Range('E45').Value = '  -6.75%  '
$ws.Range('D46').Value = '5.45'
$ws.Range('E46').Value = '  -8.05%  '
$ws.Range('E47').Value = '  -3.66%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '289.78'
$ws.Range('E48').Value = '  -7.18%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '21.33'
$ws.Range('E49').Value = '  -7.81%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  -5.46%  '
